$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "62.819.44"
$ws.Range("E2").Value = "  -0.77%  "
$ws.Range("D3").Value = "2.464.07"
$ws.Range("E3").Value = "  -0.98%  "
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "571.44"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.12%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "147.70"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.41%  "
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.531"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.68%  "
$ws.Range("E9").Value = "  -0.65%  "
$ws.Range("E10").Value = "  -0.25%  "
$ws.Range("E11").Value = "  -1.46%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.347"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.77%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "28.81"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.63%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.0000175"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.83%  "
$ws.Range("D15").Value = "2.919.85"
$ws.Range("E15").Value = "  -0.49%  "
$ws.Range("D16").Value = "62.699.50"
$ws.Range("E16").Value = "  -0.83%  "
$ws.Range("D17").Value = "2.465.37"
$ws.Range("E17").Value = "  -0.67%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "7.63"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -7.42%  "
$ws.Range("E19").Value = "  -2.96%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "321.98"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -2.45%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.14"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.20%  "
$ws.Range("E23").Value = "  -0.03%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "10.15"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +3.08%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "64.77"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -2.41%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "644.44"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -3.57%  "
$ws.Range("D27").Value = "2.595.80"
$ws.Range("E27").Value = "  -0.32%  "
$ws.Range("D28").Value = "0.0₃0963"
$ws.Range("E28").Value = "  -3.73%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.999"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.16%  "
$ws.Range("E30").Value = "  -4.27%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.87"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -2.99%  "
$ws.Range("E32").Value = "  -2.48%  "
$ws.Range("E33").Value = "  -0.58%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.999"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.04%  "
$ws.Range("E35").Value = "  -4.08%  "
$ws.Range("E36").Value = "  -2.87%  "
$ws.Range("E37").Value = "  -2.45%  "
$ws.Range("E38").Value = "  -2.06%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "150.02"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.96%  "
$ws.Range("E40").Value = "  -1.55%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.64"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -2.73%  "
$ws.Range("D44").Value = "0.0₆0305"
$ws.Range("E44").Value = "  -2.79%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "153.25"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -2.01%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "15.40"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.03%  "
$ws.Range("E47").Value = "  -1.93%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "20.28"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.41%  "
$ws.Range("E49").Value = "  -0.70%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0508"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.36%  "
$ws.Range("E51").Value = "  -1.93%  "
